$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 200000
$ws.Columns.Item(11).ColumnWidth = 7

$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("K3").Select()
